$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 250, pushing the existing rows 250-251
# (date 44544 Primera/Segunda pair) down to rows 252-253.
$ws.Rows.Item(250).Resize(2).Insert()

# New row 250: Primera, week of 44628
$ws.Cells.Item(250, 1).Value = 8
$ws.Cells.Item(250, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(250, 3).Value = "Coquimbo"
$ws.Cells.Item(250, 4).Value = 44628
$ws.Cells.Item(250, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(250, 5).Value = 4
$ws.Cells.Item(250, 6).Value = 100114014
$ws.Cells.Item(250, 7).Value = "Betarraga"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 2400
$ws.Cells.Item(250, 11).Value = 500
$ws.Cells.Item(250, 12).Value = 600
$ws.Cells.Item(250, 13).Value = 550
$ws.Cells.Item(250, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(250, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(250, 16).Value = 183
$ws.Cells.Item(250, 17).Value = 3
$ws.Cells.Item(250, 18).Value = "Hortaliza"

# New row 251: Segunda, week of 44628
$ws.Cells.Item(251, 1).Value = 8
$ws.Cells.Item(251, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(251, 3).Value = "Coquimbo"
$ws.Cells.Item(251, 4).Value = 44628
$ws.Cells.Item(251, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(251, 5).Value = 4
$ws.Cells.Item(251, 6).Value = 100114014
$ws.Cells.Item(251, 7).Value = "Betarraga"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Segunda"
$ws.Cells.Item(251, 10).Value = 1500
$ws.Cells.Item(251, 11).Value = 400
$ws.Cells.Item(251, 12).Value = 450
$ws.Cells.Item(251, 13).Value = 425
$ws.Cells.Item(251, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(251, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(251, 16).Value = 142
$ws.Cells.Item(251, 17).Value = 3
$ws.Cells.Item(251, 18).Value = "Hortaliza"
